# "Generate Report for handback" - refresh the Correspond Handoff/Handback
# datetime stamps (columns D and G) on the per-language handback-status
# sheets. Row 2 of each sheet gets a new, later pair of timestamps; row 3
# keeps its original timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-17 06:22:14"
$wsZhCn.Range("G2").Value = "2016-01-17 06:22:58"
$wsZhCn.Range("D3").Value = "2016-01-17 06:20:22"
$wsZhCn.Range("G3").Value = "2016-01-17 06:21:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-17 06:22:24"
$wsDeDe.Range("G2").Value = "2016-01-17 06:23:15"
$wsDeDe.Range("D3").Value = "2016-01-17 06:20:34"
$wsDeDe.Range("G3").Value = "2016-01-17 06:21:30"
